# Generate Report for Handback
#
# This script brings localization-status.xlsx from "Ready for handoff" to
# "Handed back" state:
#   - Overview + per-locale sheets' Status cells flip to the handback message
#   - Each locale sheet gains "Latest Target File" (F) and "Latest Handback
#     File" (G) hyperlinked entries for both data rows
#   - "Latest Handback DateTime" (H) is stamped for both rows on both locale
#     sheets (zh-cn finished a few seconds before de-de)

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both locale status columns for both rows
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value2 = $statusText
$overview.Range("C2").Value2 = $statusText
$overview.Range("B3").Value2 = $statusText
$overview.Range("C3").Value2 = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column also reflects the handback
$zhcn.Range("C2").Value2 = $statusText
$zhcn.Range("C3").Value2 = $statusText

# Row 2 (32e299a2-286b-486c-bda8-6bdc7aad0a91)
$zhcn.Hyperlinks.Add($zhcn.Cells.Item(2, 6), "https://github.com/OpenLocalizationTest/oltest/blob/9d954bef0d2f2eaa68b51e1d33fbdb50c0951d15/e2e/32e299a2-286b-486c-bda8-6bdc7aad0a91.md", "", "", "32e299a2-286b-486c-bda8-6bdc7aad0a91.md")
$zhcn.Range("F2").Font.Underline = 2
$zhcn.Range("F2").Font.Color = 15570276

$zhcn.Hyperlinks.Add($zhcn.Cells.Item(2, 7), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af5187fe2c92772389f494a4b6999a0b572b6edc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/32e299a2-286b-486c-bda8-6bdc7aad0a91.440e3be6904fdb79477d61f186bd27653b4d87db.zh-cn.xlf", "", "", "32e299a2-286b-486c-bda8-6bdc7aad0a91.440e3be6904fdb79477d61f186bd27653b4d87db.zh-cn.xlf")
$zhcn.Range("G2").Font.Underline = 2
$zhcn.Range("G2").Font.Color = 15570276

$zhcn.Range("H2").Value2 = "2016-03-24 06:52:54"

# Row 3 (968e5e94-e23d-45b5-ac89-bda0c44d0223)
$zhcn.Hyperlinks.Add($zhcn.Cells.Item(3, 6), "https://github.com/OpenLocalizationTest/oltest/blob/9d954bef0d2f2eaa68b51e1d33fbdb50c0951d15/e2e/968e5e94-e23d-45b5-ac89-bda0c44d0223.md", "", "", "968e5e94-e23d-45b5-ac89-bda0c44d0223.md")
$zhcn.Range("F3").Font.Underline = 2
$zhcn.Range("F3").Font.Color = 15570276

$zhcn.Hyperlinks.Add($zhcn.Cells.Item(3, 7), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af5187fe2c92772389f494a4b6999a0b572b6edc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/968e5e94-e23d-45b5-ac89-bda0c44d0223.4e11761e51832755bda2c16e55048fecb4012544.zh-cn.xlf", "", "", "968e5e94-e23d-45b5-ac89-bda0c44d0223.4e11761e51832755bda2c16e55048fecb4012544.zh-cn.xlf")
$zhcn.Range("G3").Font.Underline = 2
$zhcn.Range("G3").Font.Color = 15570276

$zhcn.Range("H3").Value2 = "2016-03-24 06:52:54"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Status column also reflects the handback
$dede.Range("C2").Value2 = $statusText
$dede.Range("C3").Value2 = $statusText

# Row 2 (32e299a2-286b-486c-bda8-6bdc7aad0a91)
$dede.Hyperlinks.Add($dede.Cells.Item(2, 6), "https://github.com/OpenLocalizationTest/oltest/blob/9d954bef0d2f2eaa68b51e1d33fbdb50c0951d15/e2e/32e299a2-286b-486c-bda8-6bdc7aad0a91.md", "", "", "32e299a2-286b-486c-bda8-6bdc7aad0a91.md")
$dede.Range("F2").Font.Underline = 2
$dede.Range("F2").Font.Color = 15570276

$dede.Hyperlinks.Add($dede.Cells.Item(2, 7), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e7e9a4e9a7227284d58ae19dc9835a00b1789f4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/32e299a2-286b-486c-bda8-6bdc7aad0a91.440e3be6904fdb79477d61f186bd27653b4d87db.de-de.xlf", "", "", "32e299a2-286b-486c-bda8-6bdc7aad0a91.440e3be6904fdb79477d61f186bd27653b4d87db.de-de.xlf")
$dede.Range("G2").Font.Underline = 2
$dede.Range("G2").Font.Color = 15570276

$dede.Range("H2").Value2 = "2016-03-24 06:53:03"

# Row 3 (968e5e94-e23d-45b5-ac89-bda0c44d0223)
$dede.Hyperlinks.Add($dede.Cells.Item(3, 6), "https://github.com/OpenLocalizationTest/oltest/blob/9d954bef0d2f2eaa68b51e1d33fbdb50c0951d15/e2e/968e5e94-e23d-45b5-ac89-bda0c44d0223.md", "", "", "968e5e94-e23d-45b5-ac89-bda0c44d0223.md")
$dede.Range("F3").Font.Underline = 2
$dede.Range("F3").Font.Color = 15570276

$dede.Hyperlinks.Add($dede.Cells.Item(3, 7), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4e7e9a4e9a7227284d58ae19dc9835a00b1789f4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/968e5e94-e23d-45b5-ac89-bda0c44d0223.4e11761e51832755bda2c16e55048fecb4012544.de-de.xlf", "", "", "968e5e94-e23d-45b5-ac89-bda0c44d0223.4e11761e51832755bda2c16e55048fecb4012544.de-de.xlf")
$dede.Range("G3").Font.Underline = 2
$dede.Range("G3").Font.Color = 15570276

$dede.Range("H3").Value2 = "2016-03-24 06:53:03"
